$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 18, shifting existing rows 18-34 down to 19-35
$ws.Rows.Item(18).Insert()

# Populate the newly inserted row 18 with the new weekly record
$ws.Cells.Item(18, 1).Value = 9
$ws.Cells.Item(18, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(18, 3).Value = "Metropolitana"
$ws.Cells.Item(18, 4).Value = 44438
$ws.Cells.Item(18, 5).Value = 13
$ws.Cells.Item(18, 6).Value = "Fruta"
$ws.Cells.Item(18, 7).Value = 100108
$ws.Cells.Item(18, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(18, 9).Value = 100108003
$ws.Cells.Item(18, 10).Value = "Maracuyá"
$ws.Cells.Item(18, 11).Value = "Sin especificar"
$ws.Cells.Item(18, 12).Value = "Primera"
$ws.Cells.Item(18, 13).Value = 30
$ws.Cells.Item(18, 14).Value = 32000
$ws.Cells.Item(18, 15).Value = 32000
$ws.Cells.Item(18, 16).Value = 32000
$ws.Cells.Item(18, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(18, 18).Value = "Perú"
$ws.Cells.Item(18, 19).Value = 1778
$ws.Cells.Item(18, 20).Value = 18
